$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

# Append a brand-new (still empty) paragraph right after $cur and return it.
function New-ParaAfter($cur) {
    $cur.Range.InsertParagraphAfter()
    $newp = $cur.Next()
    return $newp
}

# Add one more run of text to the end of $para (keeps current formatting of
# the range, i.e. repeated calls build up multiple runs in one paragraph).
function Add-Run($para, [string]$text) {
    $para.Range.InsertAfter($text)
}

# Create a new numbered ("Step 5" sub-list) paragraph after $cur containing
# the given run texts, and return the new paragraph.
function Add-ListPara($cur, [string[]]$runs) {
    $newp = New-ParaAfter $cur
    foreach ($t in $runs) {
        Add-Run $newp $t
    }
    $newp.Range.ListFormat.ApplyNumberDefault()
    return $newp
}

# Create a new "List Paragraph" styled paragraph (indented, no number) after
# $cur containing the given run texts, and return the new paragraph.
function Add-ListStyleOnlyPara($cur, [string[]]$runs) {
    $newp = New-ParaAfter $cur
    $newp.Range.Style = "List Paragraph"
    foreach ($t in $runs) {
        Add-Run $newp $t
    }
    return $newp
}

# ---------------------------------------------------------------------------
# Locate the insertion point: the empty paragraph that immediately precedes
# the "Launching startup home page (index)" Heading 2 paragraph.
# ---------------------------------------------------------------------------

$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r") {
        $nxt = $d.Paragraphs.Item($i + 1)
        if ($nxt.Style.NameLocal -eq "Heading 2" -and $nxt.Range.Text -eq "Launching startup home page (index)`r") {
            $anchor = $p
            break
        }
    }
}

if ($anchor -eq $null) {
    throw "Could not locate anchor paragraph before 'Launching startup home page (index)'"
}

$cur = $anchor

# ---------------------------------------------------------------------------
# Step 5 heading (bold is applied at the very end, see below, so that the
# bold paragraph-mark formatting doesn't leak into the paragraphs created
# after it).
# ---------------------------------------------------------------------------

$stepHeading = New-ParaAfter $cur
Add-Run $stepHeading "Step 5: Create Stripe "
Add-Run $stepHeading "account"
$cur = $stepHeading

# ---------------------------------------------------------------------------
# Numbered sub-steps (all share one numbered list, B..H)
# ---------------------------------------------------------------------------

$cur = Add-ListPara $cur @(
    "Go to Stripe.com and create an ",
    "account",
    " "
)

$cur = Add-ListPara $cur @(
    "After creating the ",
    "account",
    " the default mode for the account is test mode so make sure to activate the account with foundation/business information and bank details to receive real live payments"
)

$cur = Add-ListPara $cur @(
    "Create a payment link and customize you",
    "r",
    " checkout screen."
)

$cur = Add-ListPara $cur @(
    "When done, go to ",
    "product",
    " catalog and find the page you made. Go into the details and copy the ",
    "Price_ID",
    "."
)

$cur = Add-ListPara $cur @(
    "Go into the ",
    "website’s",
    " code and paste the ",
    "price_id",
    " into the ",
    "donate.php",
    " file on line 19 in between the two single quotations. Exp: ",
    "   'price' => '",
    "PASTE PRICE ID HERE",
    "'"
)

$cur = Add-ListPara $cur @(
    "Do the same for the ",
    "api",
    " secret key ",
    "under",
    " the developer section of Stripe.com."
)

$lastListPara = Add-ListPara $cur @(
    "Copy and paste the secret key into ",
    "donate.php",
    " file on line 7. "
)

# Keep building the list (J, K) contiguously first so the whole run shares
# one numId, then splice the non-numbered "Exp:" paragraph (I) in between H
# and J afterwards - this keeps J/K on the very same numbered list instead
# of starting a brand-new one.
$cur = Add-ListPara $lastListPara @(
    "Make sure to save the ",
    "donate.php",
    " file before exiting."
)

$cur = Add-ListPara $cur @(
    "If you want to change the payment link ",
    "with",
    " a different, repeat steps 3-5."
)

# Trailing blank, indented paragraph before the existing Heading 2.
$trailing = New-ParaAfter $cur
$trailing.Range.ListFormat.RemoveNumbers()
$trailing.Range.Style = "Normal"
$trailing.Range.ParagraphFormat.LeftIndent = 36

# Now insert the "Exp:" example line (List Paragraph style only, no number)
# right after "...file on line 7." and before "Make sure to save...".
Add-ListStyleOnlyPara $lastListPara @(
    "Exp:   ",
    "$",
    "stripe_secret_key",
    ' = "',
    "PASTE SECRET KEY HERE",
    '"'
) | Out-Null

# ---------------------------------------------------------------------------
# Make the "Step 5" heading bold now that every later paragraph already
# exists, so the bold paragraph-mark does not leak forward.
# ---------------------------------------------------------------------------

$stepHeading.Range.Bold = 1
$stepHeading.Range.BoldBi = 1

Write-Host "Stripe account setup section inserted; paragraphs now: $($d.Paragraphs.Count)"
